# Add the new "recipect" row (#4, dated 2019-01-04, total 320.76) to the
# sales report table on the active sheet, right below the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recipect number
$ws.Range("A5").Value = 4

# Date, stored as plain text (like the other date cells in column B) rather
# than letting Excel auto-convert the "yyyy-mm-dd" text into a date serial
# number. We build it as a text formula first, then paste-special the
# computed value back over itself so it lands as a literal shared string
# with no date number-format applied.
$ws.Range("B5").Formula = "=""2019-01-04"""
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)

# Total price
$ws.Range("C5").Value = 320.76
